$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "file / sheet / topic count" block
$ws.Range("B14").Value = "FY2019 MOC quali data"
$ws.Range("B15").Value = "Master"
$ws.Range("B16").Value = 4
$ws.Range("B17").Value = 4

# Fill in the question answers (C01..C04)
$ws.Range("B18").Value = "C01"
$ws.Range("B19").Value = "C02"
$ws.Range("B20").Value = "C03"
$ws.Range("B21").Value = "C04"

# Move the active selection to B22, as in the saved workbook
$ws.Range("B22").Select()
